$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 2) { $lastRow = 88 }

for ($r = 2; $r -le $lastRow; $r++) {
    $fVal = $ws.Cells.Item($r, 6).Value2
    if ($fVal -eq "Na") {
        $ws.Cells.Item($r, 2).Value = "No"
        $ws.Cells.Item($r, 6).Value = ""
        $ws.Cells.Item($r, 10).Value = 0
    }
}
